# Update schedule: upload lab & HW 8, shift labs 9-11 and HW 9-11 by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Week 10 area: the HW8 lecture notes duplicated into column H are no longer needed.
$ws.Range("H29").ClearContents()
$ws.Range("H30").ClearContents()

# Lec 18 gets its full title; the old duplicate in column H is removed entirely.
$ws.Range("C32").Value = "Lec 18: Outliers"
$ws.Range("H32").Clear()

# Lec 19 title shortened; HW9 link moves up from row 35 to row 33.
$ws.Range("C33").Value = "Lec 19: "

# Lab 8 (Variable Selection) is now posted, with its materials link.
$ws.Range("C31").Value = "Lab 8: Variable Selection"
$ws.Range("F31").Value = "lab-08.html"

$ws.Range("G33").Value = "hw-09"

# Former "Lab 8: Review" becomes "Lab 9: Review".
$ws.Range("C34").Value = "Lab 9: Review"

# HW9 no longer sits on the Midterm II row.
$ws.Range("G35").ClearContents()

# HW10 link moves up from row 38 to row 36.
$ws.Range("G36").Value = "hw-10"

# Former "Lab 9" becomes "Lab 10".
$ws.Range("C37").Value = "Lab 10"

# HW10 no longer sits on row 38.
$ws.Range("G38").ClearContents()

# HW11 link moves up from row 40 to row 39.
$ws.Range("G39").Value = "hw-11"

# HW11 no longer sits on row 40.
$ws.Range("G40").ClearContents()

# Former "Lab 10" becomes "Lab 11".
$ws.Range("C42").Value = "Lab 11"

# Restore the scroll position / active selection as last left by the author.
$ws.Activate()
$ws.Range("A18").Select()
$excel.ActiveWindow.ScrollRow = 18
$ws.Range("C42").Select()
